$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new date column (BR) mirroring the existing BQ ("11-sep") column,
# with header "12-sep" and the corresponding data values per row.

# Header cell - copy style from BQ1 (the previous last date column) then set value
$ws.Range("BQ1").Copy() | Out-Null
$ws.Range("BR1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("BR1").Value = "12-sep"

# Data values for rows 2-18
$ws.Range("BR2").Value = 0
$ws.Range("BR3").Value = 18.625306659136701
$ws.Range("BR4").Value = 17.413594215947633
$ws.Range("BR5").Value = 17.115094817227014
$ws.Range("BR6").Value = 0
$ws.Range("BR7").Value = 14.455851615152497
$ws.Range("BR8").Value = 13.405673647040658
$ws.Range("BR9").Value = 15.195709804162837
$ws.Range("BR10").Value = 14.668531642727368
$ws.Range("BR11").Value = 13.927406891572478
$ws.Range("BR12").Value = 0
$ws.Range("BR13").Value = 11.285700748483132
$ws.Range("BR14").Value = 0
$ws.Range("BR15").Value = 0
$ws.Range("BR16").Value = 9.597168855497312
$ws.Range("BR17").Value = 0
$ws.Range("BR18").Value = 0

# Update selection to reflect the cell that was active after the edit
$ws.Range("BT6").Select() | Out-Null

$wb.Save()
